$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the header row: Title/Author/Price -> Date/Description/Amount
# (a "Category" column label is also introduced as an available option,
# though it isn't applied to any existing column)
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Amount"

$wb.Save()
